# Updated Master data as per 16th May Refresh
# Adds three new rows (34-36) of reg_center_user_machine_h test data,
# mirroring the existing pattern (regcntr_id=10005, incrementing machine_id).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(34, 10005, 110033, 10005),
    @(35, 10005, 110034, 10005),
    @(36, 10005, 110035, 10005)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]          # A: regcntr_id
    $ws.Cells.Item($rowNum, 2).Value = $r[2]          # B: usr_id
    $ws.Cells.Item($rowNum, 3).Value = $r[3]          # C: machine_id
    $ws.Cells.Item($rowNum, 4).Value = "eng"          # D: lang_code
    $ws.Cells.Item($rowNum, 5).Value = $true          # E: is_active
    $ws.Cells.Item($rowNum, 6).Value = "superadmin"   # F: cr_by
    $ws.Cells.Item($rowNum, 7).Value = "now()"        # G: cr_dtimes
    $ws.Cells.Item($rowNum, 8).Value = "now()"        # H: eff_dtimes
}

# Reflect the post-entry selection state recorded in the workbook
# (user clicked the next empty row / selected the remainder of the sheet).
$null = $ws.Range("A37:XFD1048576").Select()
